$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center-align the header cells B2:C2 (picks up existing "horizontal center" style)
$ws.Range("B2:C2").HorizontalAlignment = -4108

# Add a new LOG row (row 18) documenting a new task
$ws.Range("A18").Value = "Iedereen"
$ws.Range("B18").Value = "Nieuwe user stories maken"
$ws.Range("C18").Value = "Eind eerste week sprint"
$ws.Range("D18").Value = "4 uur"
$ws.Range("E18").Value = "NVT"
$ws.Range("B18:E18").HorizontalAlignment = -4108

# Set the page setup (paper size A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
